$d = $word.ActiveDocument

# --- Locate the "UC007: registrer stævne" bullet (the last item in the Suggestions list) ---
$anchor = $d.Content
$found = $anchor.Find.Execute("UC007: registrer stævne", $true, $false, $false, $false, $false, `
                               $true, 1, $false, "", 0)

# --- Append two new bullets after it, in the same "Listeafsnit" list (numId 1 / ilvl 0) ---
$anchor.InsertParagraphAfter()
$firstNewIdx = $anchor.Paragraphs.Item(1).Index + 1

$para1 = $d.Paragraphs.Item($firstNewIdx)
$para1.Range.Text = "Formatér tekst i filerne"

$para1.Range.InsertParagraphAfter()
$para2 = $d.Paragraphs.Item($firstNewIdx + 1)
$para2.Range.Text = "Evt. ny UC til print stævne"

# --- The document carries a "_GoBack" bookmark right after the last bullet's text; move it
#     so it again trails the (now last) bullet instead of the original "UC007" one. ---
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

$para2 = $d.Paragraphs.Item($firstNewIdx + 1)
$insertPos = $para2.Range.End - 1

# Adding a bookmark exactly at "end of text, right before the paragraph mark" needs a
# placeholder character to anchor to, otherwise the collapsed range resolves incorrectly;
# insert one, plant the bookmark, then remove the placeholder again.
$placeholder = $d.Range($insertPos, $insertPos)
$placeholder.InsertBefore("#")
$bmRange = $d.Range($insertPos, $insertPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
$d.Range($insertPos, $insertPos + 1).Delete()
